# Update the "want to go" (想去人数) counts in column F across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets.
# "本地生活" (sheet3) has no changes.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 129
$ws1.Range("F3").Value = 1293
$ws1.Range("F5").Value = 980
$ws1.Range("F6").Value = 1748
$ws1.Range("F7").Value = 387
$ws1.Range("F8").Value = 1166
$ws1.Range("F9").Value = 51
$ws1.Range("F11").Value = 118
$ws1.Range("F12").Value = 268
$ws1.Range("F13").Value = 51
$ws1.Range("F15").Value = 650
$ws1.Range("F20").Value = 322
$ws1.Range("F21").Value = 112
$ws1.Range("F22").Value = 649
$ws1.Range("F23").Value = 18
$ws1.Range("F24").Value = 633
$ws1.Range("F25").Value = 143
$ws1.Range("F27").Value = 851
$ws1.Range("F28").Value = 303
$ws1.Range("F29").Value = 133
$ws1.Range("F30").Value = 30
$ws1.Range("F31").Value = 255
$ws1.Range("F32").Value = 8
$ws1.Range("F34").Value = 400

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 312

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 129
$ws4.Range("F4").Value = 1293
$ws4.Range("F6").Value = 980
$ws4.Range("F7").Value = 1748
$ws4.Range("F8").Value = 387
$ws4.Range("F9").Value = 1166
$ws4.Range("F10").Value = 51
$ws4.Range("F13").Value = 118
$ws4.Range("F14").Value = 268
$ws4.Range("F15").Value = 51
$ws4.Range("F17").Value = 650
$ws4.Range("F22").Value = 312
$ws4.Range("F25").Value = 322
$ws4.Range("F29").Value = 112
$ws4.Range("F30").Value = 649
$ws4.Range("F31").Value = 18
$ws4.Range("F32").Value = 633
$ws4.Range("F33").Value = 143
$ws4.Range("F35").Value = 852
$ws4.Range("F36").Value = 303
$ws4.Range("F39").Value = 133
$ws4.Range("F40").Value = 30
$ws4.Range("F41").Value = 255
$ws4.Range("F45").Value = 8
$ws4.Range("F48").Value = 400
